# Creación del constructor n
# Update the "nVeces" sample sizes (column A) and the recomputed
# backtracking times (column B, derived from the new raw millisecond
# counts over I10) on the "Backtracking" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Backtracking")

# Row 3: n = 5, t = 409 / I10
$ws.Range("A3").Value = 5
$ws.Range("B3").Formula = "=409/I10"

# Row 4: n = 10, t = 337 / I10
$ws.Range("A4").Value = 10
$ws.Range("B4").Formula = "=337/I10"

# Row 5: n = 15, t = 269 / I10
$ws.Range("A5").Value = 15
$ws.Range("B5").Formula = "=269/I10"

# Row 6: n = 25, t = 269 / I10
$ws.Range("A6").Value = 25
$ws.Range("B6").Formula = "=269/I10"

# Row 7: n = 30, t = 298 / I10 (previously an erroring formula referencing N9)
$ws.Range("A7").Value = 30
$ws.Range("B7").Formula = "=298/I10"

# Move the active selection from B3 to B7, matching the saved view state.
$ws.Range("B7").Select()

$wb.Save()
